$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# The row "VodafonePivotByCCSheetName / Pivot by CC" (row 16) is no longer
# used (the Vodafone "Pivot by CC" sheet concept was dropped), so delete the
# whole row - this shifts every row below it up by one and also drops the
# trailing blank row 1008 from the sheet's used range.
$ws.Rows.Item(16).Delete()

# A handful of the remaining "Vodafone..." prefixed setting names are
# renamed/generalised now that the template is no longer Vodafone specific:
#   VodafoneWorkingsSheetName -> InvSheetName   (value "Workings" -> "Inv")
#   VodafoneActiveListSheetName -> ActiveListSheetName
#   VodafoneWorkingsRange -> InvRange
#   VodafoneActiveListRange -> ActiveListRange
#   VodafoneActiveListReadRange -> ActiveListReadRange
$ws.Range("A14").Value = "InvSheetName"
$ws.Range("B14").Value = "Inv"

$ws.Range("A15").Value = "ActiveListSheetName"

$ws.Range("A18").Value = "InvRange"

$ws.Range("A20").Value = "ActiveListRange"

$ws.Range("A21").Value = "ActiveListReadRange"

# Restore the on-screen selection/scroll position left by the author after
# making these edits.
$ws.Range("B43").Select()
